$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asctb_incorrect_cts")

# Remove row 2, which contained all NA values; this shifts rows 3:14 up to 2:13
$ws.Rows.Item(2).Delete()
